$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 80448771
$ws.Range("B2").Value = 73693
$ws.Range("E2").Value = 6440
$ws.Range("F2").Value = "Vitgrynig nållav"
$ws.Range("G2").Value = "Chaenotheca subroscida"
$ws.Range("H2").Value = "(Eitner) Zahlbr."
$ws.Range("Q2").Value = 423289.9356373397
$ws.Range("R2").Value = 6752041.978126496
$ws.Range("AC2").ClearContents()
$ws.Range("A3").Value = 80448772
$ws.Range("B3").Value = 81236
$ws.Range("E3").Value = 1312
$ws.Range("F3").Value = "Gammelgransskål"
$ws.Range("G3").Value = "Pseudographis pinicola"
$ws.Range("H3").Value = "(Nyl.) Rehm"
$ws.Range("Q3").Value = 423289.9356373397
$ws.Range("R3").Value = 6752041.978126496
$ws.Range("AC3").ClearContents()
$ws.Range("A4").Value = 80448769
$ws.Range("Q4").Value = 422991.0759451608
$ws.Range("R4").Value = 6752021.173145968
$ws.Range("AC4").Value = "Rikligt, hkb"
$ws.Range("A5").Value = 80448775
$ws.Range("B5").Value = 77506
$ws.Range("E5").Value = 6425
$ws.Range("F5").Value = "Garnlav"
$ws.Range("G5").Value = "Alectoria sarmentosa"
$ws.Range("H5").Value = "(Ach.) Ach."
$ws.Range("Q5").Value = 423036.1594514723
$ws.Range("R5").Value = 6752009.000504656
$ws.Range("AC5").Value = "Rikligt, hkb"
$ws.Range("A6").Value = 80448777
$ws.Range("Q6").Value = 423115.1561234437
$ws.Range("R6").Value = 6752009.239606674
$ws.Range("AC6").ClearContents()
$ws.Range("A7").Value = 80448779
$ws.Range("Q7").Value = 422962.8083476268
$ws.Range("R7").Value = 6752021.785183201
$ws.Range("AC7").Value = "Rikligt"
$ws.Range("A8").Value = 80448780
$ws.Range("B8").Value = 77506
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 423056.1482692101
$ws.Range("R8").Value = 6751963.779848268
$ws.Range("AC8").Value = "Spritt"
$ws.Range("A9").Value = 80448778
$ws.Range("B9").Value = 56395
$ws.Range("C9").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("Q9").Value = 423115.1561234437
$ws.Range("R9").Value = 6752009.239606674
$ws.Range("AJ9").ClearContents()
$ws.Range("AK9").ClearContents()
$ws.Range("AO9").ClearContents()
$ws.Range("A10").Value = 80448773
$ws.Range("B10").Value = 77506
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 422635.9957601223
$ws.Range("R10").Value = 6751949.037152009
$ws.Range("AJ10").Value = "vanlig tall"
$ws.Range("AK10").Value = "Pinus sylvestris var. sylvestris"
$ws.Range("AO10").Value = "Pinus sylvestris var. sylvestris"
